$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A..K) to (B..L)
$ws.Columns.Item(1).Insert()

# Copy the style from header cell C1 onto the new header cell B1
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# New header cell for inserted column
$ws.Cells.Item(1, 2).Value = "segments"

# Copy the style from the (shifted) name column B onto the new index column A (formats only)
$ws.Range("B2:B20").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)

# Fill new column A with the numeric segment index (0-based)
for ($i = 0; $i -le 18; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}

$excel.CutCopyMode = $false

$wb.Save()
